$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3288063333333333
$ws.Range("H2").Value = 0.9864189999999999
$ws.Range("I2").Value = 0.05575527297994041
$ws.Range("J2").Value = 0.05575527297994041
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 4.709125178046333
$ws.Range("R2").Value = 42.382126602417
$ws.Range("S2").Value = 0.01644540346247579
$ws.Range("T2").Value = 0.01644540346247579
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3288063333333333
$ws.Range("H3").Value = 0.9864189999999999
$ws.Range("I3").Value = 0.05575527297994041
$ws.Range("J3").Value = 0.05575527297994041
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 8.905446300270334
$ws.Range("R3").Value = 80.149016702433
$ws.Range("S3").Value = 0.03109997120146995
$ws.Range("T3").Value = 0.03109997120146995
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3288063333333333
$ws.Range("H4").Value = 0.9864189999999999
$ws.Range("I4").Value = 0.05575527297994041
$ws.Range("J4").Value = 0.05575527297994041
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 2.350896343605444
$ws.Range("R4").Value = 21.158067092449
$ws.Range("S4").Value = 0.008209898315994663
$ws.Range("T4").Value = 0.008209898315994663
$ws.Range("I5").Value = 0.3115445049245869
$ws.Range("J5").Value = 0.3115445049245869
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 26.31324346219567
$ws.Range("R5").Value = 236.819191159761
$ws.Range("S5").Value = 0.09189220689217013
$ws.Range("T5").Value = 0.09189220689217013
$ws.Range("I6").Value = 0.3115445049245869
$ws.Range("J6").Value = 0.3115445049245869
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("R6").Value = 447.8497619002891
$ws.Range("S6").Value = 0.1737777364056091
$ws.Range("T6").Value = 0.1737777364056091
$ws.Range("I7").Value = 0.3115445049245869
$ws.Range("J7").Value = 0.3115445049245869
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 13.13613580120189
$ws.Range("R7").Value = 118.225222210817
$ws.Range("S7").Value = 0.0458745616268076
$ws.Range("T7").Value = 0.04587456162680759
$ws.Range("G8").Value = 3.731231666666667
$ws.Range("H8").Value = 11.193695
$ws.Range("I8").Value = 0.6327002220954728
$ws.Range("J8").Value = 0.6327002220954728
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 53.43825591343168
$ws.Range("R8").Value = 480.944303220885
$ws.Range("S8").Value = 0.1866193073236607
$ws.Range("T8").Value = 0.1866193073236607
$ws.Range("G9").Value = 3.731231666666667
$ws.Range("H9").Value = 11.193695
$ws.Range("I9").Value = 0.6327002220954728
$ws.Range("J9").Value = 0.6327002220954728
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 101.0573090381517
$ws.Range("R9").Value = 909.5157813433652
$ws.Range("S9").Value = 0.35291655182842
$ws.Range("T9").Value = 0.35291655182842
$ws.Range("G10").Value = 3.731231666666667
$ws.Range("H10").Value = 11.193695
$ws.Range("I10").Value = 0.6327002220954728
$ws.Range("J10").Value = 0.6327002220954728
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 26.67752410176056
$ws.Range("R10").Value = 240.097716915845
$ws.Range("S10").Value = 0.09316436294339211
$ws.Range("T10").Value = 0.0931643629433921
